# Insert a new column before the "STATUS" column (D) for the new "ETAPA" field,
# which pushes the existing STATUS column from D to E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D").Insert()
$ws.Range("D1").Value = "ETAPA"

# Match the new column's width to the narrow "ETAPA" header (closest reachable
# value to the authored 7.35 given the engine's column-width quantization).
$ws.Columns("D").ColumnWidth = 6.5

# Move the active selection, as recorded after the edit.
$ws.Range("E6").Select() | Out-Null

# Header/footer font style tweak: "Regular" -> "Normal".
$ws.PageSetup.CenterHeader = '&"Times New Roman,Normal"&12&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Normal"&12Página &P'
